$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shuffle rows: the old "Total" row (row 8) moves down to row 9, a new
# --- expense entry is written into row 7 (previously a blank spacer row),
# --- and row 8 becomes the new blank spacer row. ---

# 1) Copy the formatting (styles) of the old Total row (B8:D8) down onto the
#    new Total row (B9:D9), so row 9 ends up with the same style indices
#    (blank/bordered, bold "Total" label cell, bordered amount cell).
$ws.Range("B8:D8").Copy()
$ws.Range("B9:D9").PasteSpecial(-4122)

# 2) Copy the blank-spacer formatting from C7 onto C8, so C8 changes from the
#    bold "Total"-label style to the plain bordered style used by the spacer
#    row.
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# 3) Copy the date-number formatting from B6 onto B7, so the new entry's date
#    cell gets the date display format instead of the plain spacer style.
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# 4) Remove the old Total label/formula now that they have been relocated.
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

# 5) Write the relocated Total row contents.
$ws.Range("C9").Value = "Total"
$ws.Range("D9").Formula = "=ROUND(SUM(D3:D7),2)"

# 6) Fill in the new expense entry.
$ws.Range("B7").Value = 45301
$ws.Range("C7").Value = "Jio recharge (90737 00094) "
$ws.Range("D7").Value = 239

# 7) Widen column C slightly to fit the new, longer expense description.
$ws.Columns.Item(3).ColumnWidth = 22.65

# 8) Match the saved selection state.
$ws.Range("C8").Select()
